# 4.2.1.xlsx -- capitalize the header labels for the three "by sex" /
# "education of mother" / "wealth quintile" sub-header rows on the single
# worksheet (row 14, row 17, row 23 across columns A/B/C -- Kyrgyz/Russian/
# English). The shared-strings table in the source diff also shows these
# strings being de-duplicated/reordered, but that is an artifact of the
# authoring tool's string-table packing, not a content change: the visible
# cell text is identical before/after except for the leading capital
# letter, so we simply overwrite the affected cells with their new text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: "(жынысы/по полу/by) sex" sub-header
$ws.Range("A14").Value = "Жынысы боюнча"
$ws.Range("B14").Value = "По полу"
$ws.Range("C14").Value = "By sex"

# Row 17: "education of mother" sub-header
$ws.Range("A17").Value = "Энесинин билими "
$ws.Range("B17").Value = "Образование матери "
$ws.Range("C17").Value = "Education of mother"

# Row 23: "wealth quintile" sub-header (column A text was already
# capitalized in the source file, so only B/C change)
$ws.Range("B23").Value = "Квинтиль по индексу благосостояния"
$ws.Range("C23").Value = "Wealth quintile"

# The author's workbook was also re-saved with the cursor back on A1
# (the stray "A23" selection left over from editing is gone in the diff).
$ws.Range("A1").Select() | Out-Null
